$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.724.22"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "3.804.91"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  +3.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000254"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "4.438.66"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "3.789.01"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").Value = "67.709.77"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("E23").Value = "  +6.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "3.951.10"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +5.57%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.57%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.138"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "393.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.43%  "
